$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 5485.227
$ws.Range("I18").Value = 1880.75
$ws.Range("J18").Value = 15097.167
$ws.Range("K18").Value = 1880.75
$ws.Range("L18").Value = 15097.167
$ws.Range("M18").Value = -1596.75
$ws.Range("N18").Value = -15665.167

$ws.Range("H32").Value = 1569.7
$ws.Range("I32").Value = 2150
$ws.Range("J32").Value = 1424.625
$ws.Range("K32").Value = 2150
$ws.Range("L32").Value = 1424.625
$ws.Range("M32").Value = -1824
$ws.Range("N32").Value = -2076.625

$ws.Range("H86").Value = 103837
$ws.Range("I86").Value = 206449.17
$ws.Range("J86").Value = 1224.8334
$ws.Range("K86").Value = 206449.17
$ws.Range("L86").Value = 1224.8334
$ws.Range("M86").Value = -205326.17
$ws.Range("N86").Value = -3470.8334

$ws.Range("H88").Value = 4269.077
$ws.Range("I88").Value = 1400
$ws.Range("J88").Value = 5129.8
$ws.Range("K88").Value = 1400
$ws.Range("L88").Value = 5129.8
$ws.Range("M88").Value = -994
$ws.Range("N88").Value = -5941.8

$ws.Range("H89").Value = 103837
$ws.Range("I89").Value = 206449.17
$ws.Range("J89").Value = 1224.8334
$ws.Range("K89").Value = 1032245.85
$ws.Range("L89").Value = 6124.166999999999
$ws.Range("M89").Value = -1026629.85
$ws.Range("N89").Value = -17356.167

$ws.Range("H91").Value = 4269.077
$ws.Range("I91").Value = 1400
$ws.Range("J91").Value = 5129.8
$ws.Range("K91").Value = 1400
$ws.Range("L91").Value = 5129.8
$ws.Range("M91").Value = 4
$ws.Range("N91").Value = -7937.8

$ws.Range("H113").Value = 51417.5
$ws.Range("I113").Value = 75751.25
$ws.Range("K113").Value = 75751.25
$ws.Range("M113").Value = -72497.25

$ws.Range("H132").Value = 1572.0233
$ws.Range("I132").Value = 1473
$ws.Range("K132").Value = 4419
$ws.Range("M132").Value = -1889

$ws.Range("H135").Value = 598
$ws.Range("I135").Value = 598
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 5382
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -2847
$ws.Range("N135").Value = ""

$ws.Range("H137").Value = 1099.9791
$ws.Range("I137").Value = 893.55554
$ws.Range("K137").Value = 2680.66662
$ws.Range("M137").Value = -130.66662

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3208.7654
$ws.Range("I32").Value = 2691.282
$ws.Range("J32").Value = 16663.334
$ws.Range("K32").Value = 2691.282
$ws.Range("L32").Value = 16663.334
$ws.Range("M32").Value = -2404.282
$ws.Range("N32").Value = -17237.334

$ws.Range("H61").Value = 2723.28
$ws.Range("I61").Value = 1981
$ws.Range("J61").Value = 8166.6665
$ws.Range("K61").Value = 1981
$ws.Range("L61").Value = 8166.6665
$ws.Range("M61").Value = -1769
$ws.Range("N61").Value = -8590.666499999999

$ws.Range("H74").Value = 1775.92
$ws.Range("I74").Value = 1447
$ws.Range("K74").Value = 1447
$ws.Range("M74").Value = -573

$ws.Range("H77").Value = 1775.92
$ws.Range("I77").Value = 1447
$ws.Range("K77").Value = 7235
$ws.Range("M77").Value = -2867

$ws.Range("H82").Value = 74999.25
$ws.Range("J82").Value = 77777
$ws.Range("L82").Value = 77777
$ws.Range("N82").Value = -78499

$ws.Range("H85").Value = 74999.25
$ws.Range("J85").Value = 77777
$ws.Range("L85").Value = 77777
$ws.Range("N85").Value = -80273

$ws.Range("H136").Value = 2723.28
$ws.Range("I136").Value = 1981
$ws.Range("J136").Value = 8166.6665
$ws.Range("K136").Value = 5943
$ws.Range("L136").Value = 24499.9995
$ws.Range("M136").Value = -3393
$ws.Range("N136").Value = -29599.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 85000
$ws.Range("J108").Value = 85000
$ws.Range("L108").Value = 85000
$ws.Range("N108").Value = -92680

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1748.7941
$ws.Range("I31").Value = 1355.1305
$ws.Range("J31").Value = 2571.9092
$ws.Range("K31").Value = 1355.1305
$ws.Range("L31").Value = 2571.9092
$ws.Range("M31").Value = -1060.1305
$ws.Range("N31").Value = -3161.9092

$ws.Range("H34").Value = 1748.7941
$ws.Range("I34").Value = 1355.1305
$ws.Range("J34").Value = 2571.9092
$ws.Range("K34").Value = 1355.1305
$ws.Range("L34").Value = 2571.9092
$ws.Range("M34").Value = -1153.1305
$ws.Range("N34").Value = -2975.9092

$ws.Range("H58").Value = 1978043.9
$ws.Range("I58").Value = 3345646.8
$ws.Range("J58").Value = 2617.4443
$ws.Range("K58").Value = 3345646.8
$ws.Range("L58").Value = 2617.4443
$ws.Range("M58").Value = -3345443.8
$ws.Range("N58").Value = -3023.4443

$ws.Range("H99").Value = 1747
$ws.Range("I99").Value = 1747
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1747
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -249
$ws.Range("N99").Value = ""

$ws.Range("H126").Value = 1747
$ws.Range("I126").Value = 1747
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5241
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -2771
$ws.Range("N126").Value = ""

$ws.Range("H134").Value = 1728.5
$ws.Range("I134").Value = 1617.8918
$ws.Range("K134").Value = 4853.6754
$ws.Range("M134").Value = -2318.6754

$ws.Range("H136").Value = 1978043.9
$ws.Range("I136").Value = 3345646.8
$ws.Range("J136").Value = 2617.4443
$ws.Range("K136").Value = 10036940.4
$ws.Range("L136").Value = 7852.3329
$ws.Range("M136").Value = -10034390.4
$ws.Range("N136").Value = -12952.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").Value = ""

$ws.Range("H100").Value = 3151.6667
$ws.Range("I100").Value = 2025
$ws.Range("J100").Value = 3377
$ws.Range("K100").Value = 6075
$ws.Range("L100").Value = 10131
$ws.Range("M100").Value = -5264
$ws.Range("N100").Value = -11753

$ws.Range("H103").Value = 2310.6667
$ws.Range("I103").Value = 2462
$ws.Range("J103").Value = 2189.6
$ws.Range("K103").Value = 7386
$ws.Range("L103").Value = 6568.799999999999
$ws.Range("M103").Value = -6507
$ws.Range("N103").Value = -8326.799999999999

$ws.Range("H106").Value = 8000
$ws.Range("J106").Value = 8000
$ws.Range("L106").Value = 24000
$ws.Range("N106").Value = -25892

$ws.Range("H132").Value = 943.5
$ws.Range("I132").Value = 849.6667
$ws.Range("K132").Value = 7647.0003
$ws.Range("M132").Value = -5117.0003

$ws.Range("H136").Value = 1692.4546
$ws.Range("I136").Value = 1692.4546
$ws.Range("K136").Value = 5077.3638
$ws.Range("M136").Value = 22.63619999999992

$ws.Range("H138").Value = 2371.5557
$ws.Range("I138").Value = 2134.6875
$ws.Range("K138").Value = 6404.0625
$ws.Range("M138").Value = -1264.0625

$ws.Range("H139").Value = 5767.96
$ws.Range("I139").Value = 5924.9585
$ws.Range("J139").Value = 2000
$ws.Range("K139").Value = 17774.8755
$ws.Range("L139").Value = 6000
$ws.Range("M139").Value = -12634.8755
$ws.Range("N139").Value = -16280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2573283.8
$ws.Range("I126").Value = 13892450
$ws.Range("K126").Value = 41677350
$ws.Range("M126").Value = -41674880

$ws.Range("H136").Value = 8940.809999999999
$ws.Range("J136").Value = 8940.809999999999
$ws.Range("L136").Value = 26822.43
$ws.Range("N136").Value = -31922.43

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7867.095
$ws.Range("I40").Value = 8082.5
$ws.Range("K40").Value = 8082.5
$ws.Range("M40").Value = -7946.5

$ws.Range("H46").Value = 1889.5385
$ws.Range("J46").Value = 2845.3333
$ws.Range("L46").Value = 2845.3333
$ws.Range("N46").Value = -3221.3333

$ws.Range("H61").Value = 2458.95
$ws.Range("I61").Value = 2366.25
$ws.Range("K61").Value = 2366.25
$ws.Range("M61").Value = -2164.25

$ws.Range("H98").Value = 100000
$ws.Range("J98").Value = 100000
$ws.Range("L98").Value = 100000
$ws.Range("N98").Value = -105990

$ws.Range("H113").Value = 2458.95
$ws.Range("I113").Value = 2366.25
$ws.Range("K113").Value = 2366.25
$ws.Range("M113").Value = -196.25

$ws.Range("H122").Value = 2757.5293
$ws.Range("I122").Value = 1717.8
$ws.Range("J122").Value = 4242.857
$ws.Range("K122").Value = 5153.4
$ws.Range("L122").Value = 12728.571
$ws.Range("M122").Value = -2703.4
$ws.Range("N122").Value = -17628.571

$ws.Range("H136").Value = 2609.7441
$ws.Range("I136").Value = 1673.2188
$ws.Range("J136").Value = 5334.1816
$ws.Range("K136").Value = 5019.6564
$ws.Range("L136").Value = 16002.5448
$ws.Range("M136").Value = -2469.6564
$ws.Range("N136").Value = -21102.5448

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 19999.5
$ws.Range("J92").Value = 19999.5
$ws.Range("L92").Value = 19999.5
$ws.Range("N92").Value = -24991.5

$ws.Range("H126").Value = 3733.4
$ws.Range("I126").Value = 3663.25
$ws.Range("J126").Value = 3758.9092
$ws.Range("K126").Value = 10989.75
$ws.Range("L126").Value = 11276.7276
$ws.Range("M126").Value = -8519.75
$ws.Range("N126").Value = -16216.7276
